$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: DQ_Metrics
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("DQ_Metrics")

# Swap the "Analyzed Patients" / "Analyzed Cases" rows (aCase now comes first)
$ws1.Cells.Item(11,1).Value = "aCase"
$ws1.Cells.Item(11,2).Value = "Analyzed Cases"
$ws1.Cells.Item(12,1).Value = "aPatient"
$ws1.Cells.Item(12,2).Value = "Analyzed Patients"

# Parallel computing optimization: execution time improved, new CPU core count metric
$ws1.Cells.Item(27,3).Value = "'0.04"
$ws1.Cells.Item(27,3).Style = "Normal"

# Insert a new row for the CPU core metric right after Execution Time
$ws1.Rows.Item(28).Insert()
$ws1.Cells.Item(28,1).Value = "cpu_core"
$ws1.Cells.Item(28,2).Value = "CPU cores"
$ws1.Cells.Item(28,3).Value = "'4"
$ws1.Cells.Item(28,3).Style = "Normal"

# ---------------------------------------------------------------
# Sheet 2: DQ_Violations (refreshed / re-ordered data export)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DQ_Violations")

$violations = @(
    @("260123451-Airolo-P-0000247","260123451-Airolo-F-0000255","E84.0","","Missing Orpha Code.  "),
    @("260123451-Airolo-P-0000247","260123451-Airolo-F-0000255","E84.1","","Missing Orpha Code.  "),
    @("260123451-Airolo-P-0000248","260123451-Airolo-F-0000256","E84.1","","Missing Orpha Code.  "),
    @("260123451-Airolo-P-0000251","260123451-Airolo-F-0000259","E84.9","","Missing Orpha Code.  "),
    @("260123451-Airolo-P-0000527","260123451-Airolo-F-0000545","E84.80","","Missing Orpha Code.  "),
    @("260123451-Airolo-P-0000002","260123451-Airolo-F-0000003","D86.1","797","Implausible birthdate 1877-12-01 maximal age 130. "),
    @("260123451-Airolo-P-0000100","260123451-Airolo-F-0000104","E03.0","797","Ambiguous Orphacoding. ICD10-Orpha combination: E03.0 - 797 is implausible according to Alpha-ID-SE. "),
    @("260123451-Airolo-P-0000101","260123451-Airolo-F-0000105","E03.1","797","Ambiguous Orphacoding. ICD10-Orpha combination: E03.1 - 797 is implausible according to Alpha-ID-SE. "),
    @("260123451-Airolo-P-0000345","260123451-Airolo-F-0000354","E03.1","442","Ambiguous Orphacoding. ICD10-Orpha combination: E03.1 - 442 is implausible according to Alpha-ID-SE. ")
)

for ($i = 0; $i -lt $violations.Length; $i++) {
    $r = $i + 2
    $row = $violations[$i]
    $ws2.Cells.Item($r,1).Value = $row[0]
    $ws2.Cells.Item($r,2).Value = $row[1]
    $ws2.Cells.Item($r,3).Value = $row[2]
    if ($row[3] -eq "") {
        $ws2.Cells.Item($r,4).Value = ""
    } else {
        $ws2.Cells.Item($r,4).Value = "'" + $row[3]
        $ws2.Cells.Item($r,4).Style = "Normal"
    }
    $ws2.Cells.Item($r,5).Value = $row[4]
}
